# #327 Ajout des profils d'acces a58d18c1e8091c98efec92c8c093b361a253eee5
#
# 1. Update the "Date" metadata value on the "Metadata" sheet
#    (2024-03-14T13:39:21+00:00 -> 2024-03-19T13:17:15+00:00).
# 2. On the "Elements" sheet, swap the two "Mapping" columns (AK <-> AL):
#    header text, column widths and the per-row data all move together.

$wb = $excel.ActiveWorkbook

# --- 1. Metadata!B8 : Date value -----------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# --- 2. Elements sheet: swap columns AK (37) and AL (38) -----------------
$ws = $wb.Worksheets.Item("Elements")

$colAK = 37
$colAL = 38

# Swap the column widths so the wider "Spécification" column keeps its
# width regardless of which letter it now lives under (AK was 24.98, AL was
# 73.90 - after the swap AK becomes the wide one and AL the narrow one).
# (ColumnWidth is quantized to whole pixels by the host, so these inputs are
# the values that land closest to the exact target widths of 73.8984375 /
# 24.98046875 characters.)
$ws.Columns.Item($colAK).ColumnWidth = 73.0
$ws.Columns.Item($colAL).ColumnWidth = 24.166666666666668

# Swap the cell values in columns AK and AL, row by row (including the
# header in row 1). Rows where both columns already hold the same value are
# left completely untouched, so their underlying representation stays
# pristine (a swap there is a visual no-op in the source diff too).
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cellAK = $ws.Cells.Item($r, $colAK)
    $cellAL = $ws.Cells.Item($r, $colAL)
    $valAK = $cellAK.Value2
    $valAL = $cellAL.Value2
    if ($valAK -ne $valAL) {
        $cellAK.Value = $valAL
        $cellAL.Value = $valAK
    }
}
